$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of B2 and B3
$b2 = $ws.Range("B2").Value2
$b3 = $ws.Range("B3").Value2
$ws.Range("B2").Value = $b3
$ws.Range("B3").Value = $b2

# Move active cell selection from D12 to D13
$ws.Range("D13").Select()
